# The commit reorders the weekly price records (rows 2-33, columns A:T)
# of the "Fruta, Terminal La Palmera de La Serena - Tuna" sheet.
# Every data row's full contents move to a different row; the header
# row (row 1) is untouched. Read every source row into memory first,
# then write each one back out at its new location so that rows which
# are both a source and a destination are handled safely regardless of
# write order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current contents (A:T) of every data row before
# overwriting anything.
$orig = @{}
for ($r = 2; $r -le 33; $r++) {
    $orig[$r] = $ws.Range("A$r`:T$r").Value2
}

# Maps each destination row to the source row whose data it should
# receive, i.e. new row -> old row.
$mapping = @{
    2  = 8
    3  = 9
    4  = 18
    5  = 19
    6  = 20
    7  = 23
    8  = 24
    9  = 25
    10 = 13
    11 = 14
    12 = 17
    13 = 21
    14 = 22
    15 = 29
    16 = 30
    17 = 27
    18 = 28
    19 = 15
    20 = 16
    21 = 10
    22 = 11
    23 = 12
    24 = 4
    25 = 5
    26 = 2
    27 = 3
    28 = 31
    29 = 32
    30 = 33
    31 = 6
    32 = 7
    33 = 26
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Range("A$newRow`:T$newRow").Value2 = $orig[$oldRow]
}
